$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.321.96"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "1.683.18"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5536"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +8.27%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2703"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06502"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07561"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.545"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").Value = "1.678.52"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5811"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008462"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.17"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").Value = "26.395.54"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.939"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.91"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.21"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.49"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1328"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +10.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.896"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.58%  "
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06366"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.394"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.11%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.594"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.583"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.669"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.041"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6226"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.719"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.238"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").Value = "1.112.52"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01627"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8750"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.33%  "
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "1.833.07"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000110"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.34"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.208"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4294"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.078"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.06%  "
